$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D3").Value = "2016-01-08 09:14:44"
$zhcn.Range("G3").Value = "2016-01-08 09:15:27"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D3").Value = "2016-01-08 09:14:54"
$dede.Range("G3").Value = "2016-01-08 09:15:47"
